$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update selection on "AddNewAgent" (it is the currently active sheet)
#    before we insert the new sheet in front of it, so it stops being the
#    active tab once the new sheet is added further below.
# ---------------------------------------------------------------------------
$addNewAgent = $wb.Worksheets.Item("AddNewAgent")
$addNewAgent.Range("D23").Select()

# ---------------------------------------------------------------------------
# 2) Update selection on "CallCenter" (no activation needed for the write,
#    selecting is only way to change stored sqref/activeCell).
# ---------------------------------------------------------------------------
$callCenter = $wb.Worksheets.Item("CallCenter")
$callCenter.Range("A1:D2").Select()

# ---------------------------------------------------------------------------
# 3) Insert the new "AddAgency" worksheet right after "CallCenter". Adding it
#    last makes it the active / selected tab, matching the target workbook.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $callCenter)
$newSheet.Name = "AddAgency"

# Row 1 headers - reuse the same cell formatting already used on CallCenter.
$callCenter.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$newSheet.Range("A1").Value = "project"

$callCenter.Range("B1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$newSheet.Range("B1").Value = "TestScenario"

$callCenter.Range("C1").Copy()
$newSheet.Range("C1").PasteSpecial(-4122)
$newSheet.Range("C1").Value = "Run"

$callCenter.Range("D1").Copy()
$newSheet.Range("D1").PasteSpecial(-4122)
$newSheet.Range("D1").Value = "Zone"

$newSheet.Range("E1").Value = "Region"

# Row 2 data.
$callCenter.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A2").Value = "Beacon FCM"

$callCenter.Range("B2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)
$newSheet.Range("B2").Value = "CoreAddAgency"

$callCenter.Range("C2").Copy()
$newSheet.Range("C2").PasteSpecial(-4122)
$newSheet.Range("C2").Value = "Yes"

$newSheet.Range("D2").Value = "Ahmedabad"
$newSheet.Range("D2").Font.Name = "Courier New"
$newSheet.Range("D2").Font.Size = 9
$newSheet.Range("D2").Font.Color = 2039583

$newSheet.Range("D2").Copy()
$newSheet.Range("E2").PasteSpecial(-4122)
$newSheet.Range("E2").Value = "Indore"

$newSheet.Rows.Item(2).RowHeight = 30

$newSheet.Range("D2").Select()
